$d = $word.ActiveDocument

# Locate the paragraph that currently holds the split "Assume station..." text
# (it has the bookmark `_GoBack` sitting mid-run). We replace that paragraph's
# whole range with the corrected paragraph plus five new red "storyline" notes
# that follow it, using InsertXML so we get full control of run/bookmark layout.

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Assume station gets restocked every day*") {
        $target = $p
        break
    }
}

$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>
<w:r><w:t xml:space="preserve">Assume station gets restocked every day (make up a cost to restock station ($150)) compare it to the </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>one time</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t xml:space="preserve"> fixed cost of adding extra docks (lets say $3000) and look at monthly cost differential to implement fix</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr>
<w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Check correlation between wait time and inflow/outflow</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr>
<w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>In prediction model, double capacity and determine the number of 1&#8217;s that have been decreased or increased</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr>
<w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Run clusters with prediction</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr>
<w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Add features (</w:t></w:r>
<w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">average trip duration per start station, average trip duration for end station) </w:t></w:r>
<w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>after train test split</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
<w:p>
<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr>
<w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Make train/test split into two CSV files for everyone to work off of</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$target.Range.InsertXML($xml)
Write-Output "applied"
